$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New note row (row 26)
$ws.Cells.Item(26, 1).Value = "ab hier mit cuml"

# New training run row (row 27)
$ws.Cells.Item(27, 1).Value = "regular"
$ws.Cells.Item(27, 2).Value = "full random"
$ws.Cells.Item(27, 3).Value = "sentences"
$ws.Cells.Item(27, 4).Value = 5000
$ws.Cells.Item(27, 5).Value = 200
$ws.Cells.Item(27, 6).Value = 5
$ws.Cells.Item(27, 8).Value = "2900s"
$ws.Cells.Item(27, 9).Value = "null"
$ws.Cells.Item(27, 10).Value = "yes"
$ws.Cells.Item(27, 11).Value = 100

# Update view: zoom level and active selection (also drops stale topLeftCell)
$excel.ActiveWindow.Zoom = 76
$ws.Range("L21").Select()
